$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The 16 "mac_i2s_send_count_0".."mac_i2s_send_count_15" register rows
# (rows 21-36) are being removed. Deleting them shifts the remaining
# "mac_test_array" row (old row 37) up to row 21, and Excel automatically
# updates the sheet dimension and re-flows the row numbering below the
# deleted range.
$ws.Rows("21:36").Delete()
